$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (material/run names) from row 4 down to row 29 shift down by two
# conceptually: "Holden" and "Rizzie Spiral" are new entries inserted near the
# top of the list, pushing every later name down by two rows; "Thomas Hex" is
# also renamed to "Matthies Hex" along the way. The "A" index column and the
# 1s filling columns C:T are left untouched for these rows.
$names = @(
  "Holden",
  "Rizzie Spiral",
  "RotRing OmegaMax-90",
  "Equal Angle",
  "Tilt Rotate",
  "CLR",
  "Rizzie Hex",
  "Matthies Hex",
  "Tilt Rotate_Partial",
  "RotRing OmegaMax-60",
  "Equal Angle_Partial",
  "Rizzie Hex_Partial",
  "ND Single",
  "RD Single",
  "TD Single",
  "Morris Single",
  "Ring Perpendicular to ND",
  "Ring Perpendicular to RD",
  "Ring Perpendicular to TD",
  "OffsetFTD",
  "OffsetATD",
  "OffsetF45",
  "OffsetA45",
  "OffsetFRD",
  "OffsetARD",
  "Gaussian Quadrature"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

# Append two brand-new simulation rows (30 and 31) continuing the pattern,
# cloning row 29's formatting/values then overwriting the index + name cells.
$ws.Range("A29:T29").Copy($ws.Range("A30:T30"))
$ws.Range("A29:T29").Copy($ws.Range("A31:T31"))

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "Michael-CCHex"

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "Michael-SNHex"

Write-Host "edit applied"
